$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.NumberFormat = "General"
}

function Set-PlainCell($range, $value) {
    $range.Value = $value
}

Set-TextCell $ws.Range("D2") "69.350.26"
Set-PlainCell $ws.Range("E2") "  +0.90%  "

Set-TextCell $ws.Range("D3") "3.790.29"
Set-PlainCell $ws.Range("E3") "  +1.18%  "

Set-PlainCell $ws.Range("E4") "  +0.04%  "

Set-TextCell $ws.Range("D5") "603.85"
Set-PlainCell $ws.Range("E5") "  +0.38%  "

Set-TextCell $ws.Range("D6") "165.52"
Set-PlainCell $ws.Range("E6") "  -2.14%  "

Set-TextCell $ws.Range("D7") "3.783.77"
Set-PlainCell $ws.Range("E7") "  +1.05%  "

Set-PlainCell $ws.Range("E8") "  -0.05%  "

Set-PlainCell $ws.Range("E9") "  +0.68%  "

Set-PlainCell $ws.Range("E10") "  +4.60%  "

Set-TextCell $ws.Range("D11") "6.34"
Set-PlainCell $ws.Range("E11") "  -0.02%  "

Set-TextCell $ws.Range("D12") "0.461"
Set-PlainCell $ws.Range("E12") "  -0.27%  "

Set-TextCell $ws.Range("D13") "37.70"
Set-PlainCell $ws.Range("E13") "  -1.63%  "

Set-TextCell $ws.Range("D14") "0.0000248"
Set-PlainCell $ws.Range("E14") "  +0.14%  "

Set-TextCell $ws.Range("D15") "4.424.71"
Set-PlainCell $ws.Range("E15") "  +1.16%  "

Set-TextCell $ws.Range("D16") "3.792.73"
Set-PlainCell $ws.Range("E16") "  +1.27%  "

Set-TextCell $ws.Range("D17") "69.462.82"
Set-PlainCell $ws.Range("E17") "  +1.04%  "

Set-TextCell $ws.Range("D18") "7.44"
Set-PlainCell $ws.Range("E18") "  +2.18%  "

Set-TextCell $ws.Range("D19") "17.64"
Set-PlainCell $ws.Range("E19") "  +3.23%  "

Set-PlainCell $ws.Range("E20") "  -0.85%  "

Set-TextCell $ws.Range("D21") "11.40"
Set-PlainCell $ws.Range("E21") "  +5.63%  "

Set-TextCell $ws.Range("D22") "494.38"
Set-PlainCell $ws.Range("E22") "  -0.30%  "

Set-TextCell $ws.Range("D23") "0.726"
Set-PlainCell $ws.Range("E23") "  -0.39%  "

Set-PlainCell $ws.Range("E24") "  -1.74%  "

Set-TextCell $ws.Range("D25") "84.93"
Set-PlainCell $ws.Range("E25") "  -0.48%  "

Set-TextCell $ws.Range("D26") "2.27"
Set-PlainCell $ws.Range("E26") "  -2.43%  "

Set-PlainCell $ws.Range("E27") "  -0.09%  "

Set-TextCell $ws.Range("D28") "10.15"
Set-PlainCell $ws.Range("E28") "  -1.92%  "

Set-PlainCell $ws.Range("E29") "  +0.09%  "

Set-PlainCell $ws.Range("E30") "  -0.08%  "

Set-TextCell $ws.Range("D31") "8.13"
Set-PlainCell $ws.Range("E31") "  +2.52%  "

Set-TextCell $ws.Range("D32") "2.43"
Set-PlainCell $ws.Range("E32") "  -3.39%  "

Set-TextCell $ws.Range("D33") "31.97"
Set-PlainCell $ws.Range("E33") "  +0.39%  "

Set-TextCell $ws.Range("D34") "3.932.77"
Set-PlainCell $ws.Range("E34") "  +1.02%  "

Set-TextCell $ws.Range("D35") "3.737.20"
Set-PlainCell $ws.Range("E35") "  +1.52%  "

Set-PlainCell $ws.Range("E36") "  -0.68%  "

Set-TextCell $ws.Range("D37") "5.98"
Set-PlainCell $ws.Range("E37") "  +2.07%  "

Set-PlainCell $ws.Range("E38") "  +0.41%  "

Set-PlainCell $ws.Range("E39") "  +4.93%  "

Set-TextCell $ws.Range("D40") "1.00"
Set-PlainCell $ws.Range("E40") "  +0.04%  "

Set-TextCell $ws.Range("D41") "0.326"
Set-PlainCell $ws.Range("E41") "  +0.72%  "

Set-TextCell $ws.Range("D42") "3.08"
Set-PlainCell $ws.Range("E42") "  +5.15%  "

Set-PlainCell $ws.Range("E43") "  +1.04%  "

Set-PlainCell $ws.Range("B44") "OKB"
Set-PlainCell $ws.Range("C44") "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextCell $ws.Range("D44") "48.50"
Set-PlainCell $ws.Range("E44") "  -0.72%  "

Set-PlainCell $ws.Range("B45") "Bittensor"
Set-PlainCell $ws.Range("C45") "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextCell $ws.Range("D45") "425.55"
Set-PlainCell $ws.Range("E45") "  -2.93%  "

Set-TextCell $ws.Range("D46") "8.46"
Set-PlainCell $ws.Range("E46") "  -0.08%  "

Set-TextCell $ws.Range("D48") "40.36"
Set-PlainCell $ws.Range("E48") "  -0.62%  "

Set-PlainCell $ws.Range("B49") "Monero"
Set-PlainCell $ws.Range("C49") "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextCell $ws.Range("D49") "141.91"
Set-PlainCell $ws.Range("E49") "  +0.26%  "

Set-PlainCell $ws.Range("B50") "Maker"
Set-PlainCell $ws.Range("C50") "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextCell $ws.Range("D50") "2.816.36"
Set-PlainCell $ws.Range("E50") "  +0.96%  "

Set-TextCell $ws.Range("D51") "1.29"
Set-PlainCell $ws.Range("E51") "  +7.38%  "
